$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.811.20'
$ws.Range("E2").Value = '  +1.34%  '

$ws.Range("D3").Value = '2.087.06'
$ws.Range("E3").Value = '  +1.06%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = '234.41'
$ws.Range("E5").Value = '  -0.28%  '

$ws.Range("D6").Value = '0.625'
$ws.Range("E6").Value = '  -0.09%  '

$ws.Range("D7").Value = '58.84'
$ws.Range("E7").Value = '  +2.27%  '

$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("D9").Value = '0.391'
$ws.Range("E9").Value = '  -1.31%  '

$ws.Range("D10").Value = '0.0793'
$ws.Range("E10").Value = '  +2.44%  '

$ws.Range("E11").Value = '  +3.23%  '

$ws.Range("D12").Value = '2.394.20'
$ws.Range("E12").Value = '  +1.03%  '

$ws.Range("D13").Value = '14.72'
$ws.Range("E13").Value = '  +2.31%  '

$ws.Range("D14").Value = '21.33'
$ws.Range("E14").Value = '  +2.63%  '

$ws.Range("D15").Value = '0.770'
$ws.Range("E15").Value = '  -0.70%  '

$ws.Range("D16").Value = '5.30'
$ws.Range("E16").Value = '  +2.52%  '

$ws.Range("D17").Value = '2.092.72'
$ws.Range("E17").Value = '  +1.31%  '

$ws.Range("D18").Value = '37.712.06'
$ws.Range("E18").Value = '  +1.08%  '

$ws.Range("D19").Value = '6.26'
$ws.Range("E19").Value = '  +0.80%  '

$ws.Range("D20").Value = '71.42'
$ws.Range("E20").Value = '  +2.59%  '

$ws.Range("D21").Value = '0.0₃0830'
$ws.Range("E21").Value = '  +1.51%  '

$ws.Range("D22").Value = '228.39'
$ws.Range("E22").Value = '  +0.68%  '

$ws.Range("E23").Value = '  +0.04%  '

$ws.Range("D24").Value = '2.42'
$ws.Range("E24").Value = '  -0.70%  '

$ws.Range("D25").Value = '2.38'
$ws.Range("E25").Value = '  -0.73%  '

$ws.Range("D26").Value = '170.41'
$ws.Range("E26").Value = '  +1.97%  '

$ws.Range("E27").Value = '  +9.17%  '

$ws.Range("D28").Value = '9.07'
$ws.Range("E28").Value = '  +2.09%  '

$ws.Range("E29").Value = '  -0.15%  '

$ws.Range("D30").Value = '19.54'
$ws.Range("E30").Value = '  +2.09%  '

$ws.Range("E31").Value = '  +2.27%  '

$ws.Range("E32").Value = '  +3.96%  '

$ws.Range("E33").Value = '  +2.67%  '

$ws.Range("D34").Value = '4.70'
$ws.Range("E34").Value = '  +3.55%  '

$ws.Range("E35").Value = '  +0.82%  '

$ws.Range("D36").Value = '3.49'
$ws.Range("E36").Value = '  +5.00%  '

$ws.Range("E37").Value = '  +2.64%  '

$ws.Range("E38").Value = '  -0.10%  '

$ws.Range("D39").Value = '5.44'
$ws.Range("E39").Value = '  -3.58%  '

$ws.Range("E40").Value = '  +2.42%  '

$ws.Range("D41").Value = '99.11'
$ws.Range("E41").Value = '  +1.41%  '

$ws.Range("E42").Value = '  -0.05%  '

$ws.Range("D43").Value = '0.0215'
$ws.Range("E43").Value = '  +1.27%  '

$ws.Range("D44").Value = '1.460.85'
$ws.Range("E44").Value = '  -1.42%  '

$ws.Range("D45").Value = '1.17'
$ws.Range("E45").Value = '  +0.81%  '

$ws.Range("D46").Value = '4.32'
$ws.Range("E46").Value = '  +6.81%  '

$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D47").Value = '1.08'
$ws.Range("E47").Value = '  +5.20%  '

$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").Value = '16.10'
$ws.Range("E48").Value = '  +5.10%  '

$ws.Range("D49").Value = '7.48'
$ws.Range("E49").Value = '  +3.40%  '

$ws.Range("E50").Value = '  +2.86%  '

$ws.Range("D51").Value = '47.39'
$ws.Range("E51").Value = '  +5.67%  '
